# Update the "data" sheet F column (time_taken) timestamps for rows 2-70
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$timestamps = @(
    "2021-10-05 14:22:08.155107",
    "2021-10-05 14:22:08.155114",
    "2021-10-05 14:22:08.155117",
    "2021-10-05 14:22:08.155120",
    "2021-10-05 14:22:08.155123",
    "2021-10-05 14:22:08.155125",
    "2021-10-05 14:22:08.155128",
    "2021-10-05 14:22:08.155130",
    "2021-10-05 14:22:08.155133",
    "2021-10-05 14:22:08.155136",
    "2021-10-05 14:22:08.155138",
    "2021-10-05 14:22:08.155141",
    "2021-10-05 14:22:08.155144",
    "2021-10-05 14:22:08.155146",
    "2021-10-05 14:22:08.155148",
    "2021-10-05 14:22:08.155151",
    "2021-10-05 14:22:08.155154",
    "2021-10-05 14:22:08.155156",
    "2021-10-05 14:22:08.155159",
    "2021-10-05 14:22:08.155161",
    "2021-10-05 14:22:08.155164",
    "2021-10-05 14:22:08.155166",
    "2021-10-05 14:22:08.155169",
    "2021-10-05 14:22:08.155171",
    "2021-10-05 14:22:08.155174",
    "2021-10-05 14:22:08.155176",
    "2021-10-05 14:22:08.155179",
    "2021-10-05 14:22:08.155181",
    "2021-10-05 14:22:08.155184",
    "2021-10-05 14:22:08.155186",
    "2021-10-05 14:22:08.155188",
    "2021-10-05 14:22:08.155191",
    "2021-10-05 14:22:08.155194",
    "2021-10-05 14:22:08.155196",
    "2021-10-05 14:22:08.155199",
    "2021-10-05 14:22:08.155201",
    "2021-10-05 14:22:08.155204",
    "2021-10-05 14:22:08.155206",
    "2021-10-05 14:22:08.155209",
    "2021-10-05 14:22:08.155211",
    "2021-10-05 14:22:08.155214",
    "2021-10-05 14:22:08.155216",
    "2021-10-05 14:22:08.155219",
    "2021-10-05 14:22:08.155221",
    "2021-10-05 14:22:08.155224",
    "2021-10-05 14:22:08.155226",
    "2021-10-05 14:22:08.155229",
    "2021-10-05 14:22:08.155231",
    "2021-10-05 14:22:08.155234",
    "2021-10-05 14:22:08.155236",
    "2021-10-05 14:22:08.155238",
    "2021-10-05 14:22:08.155241",
    "2021-10-05 14:22:08.155244",
    "2021-10-05 14:22:08.155246",
    "2021-10-05 14:22:08.155249",
    "2021-10-05 14:22:08.155251",
    "2021-10-05 14:22:08.155254",
    "2021-10-05 14:22:08.155256",
    "2021-10-05 14:22:08.155259",
    "2021-10-05 14:22:08.155261",
    "2021-10-05 14:22:08.155264",
    "2021-10-05 14:22:08.155266",
    "2021-10-05 14:22:08.155268",
    "2021-10-05 14:22:08.155271",
    "2021-10-05 14:22:08.155274",
    "2021-10-05 14:22:08.155277",
    "2021-10-05 14:22:08.155280",
    "2021-10-05 14:22:08.155282",
    "2021-10-05 14:22:08.155285"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

# Add the new "metadata" sheet after "data"
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$metaSheet.Name = "metadata"

# Header row values
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Copy the header formatting (bold, bordered, centered) from the "data" sheet header row
$ws.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)

# Data row
$metaSheet.Range("A2").Value = 0
$ws.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$metaSheet.Range("B2").Value = "Palmoplantar keratodermas"
$metaSheet.Range("C2").Value = 556
$metaSheet.Range("D2").Value = "'1.9"
$metaSheet.Range("D2").Style = "Normal"
$metaSheet.Range("E2").Value = "2021-08-31T13:41:27.204199Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:22:08.151668"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/556/?format=json"

# Keep "data" as the active/selected sheet (unchanged in the target workbook view)
$ws.Activate()

Write-Output "edit complete"
